$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# ALC row 4 (G=5470)
$wsALC.Range("H4").Value = 226.88889
$wsALC.Range("I4").Value = 178.33333
$wsALC.Range("J4").Value = 324
$wsALC.Range("K4").Value = 178.33333
$wsALC.Range("L4").Value = 324
$wsALC.Range("M4").Value = -64.33332999999999
$wsALC.Range("N4").Value = -552

# ALC row 10 (G=1959)
$wsALC.Range("H10").Value = 31354.889
$wsALC.Range("J10").Value = 35248.75
$wsALC.Range("L10").Value = 35248.75
$wsALC.Range("N10").Value = -35834.75

# ALC row 20 (G=1965)
$wsALC.Range("H20").Value = 1294.4286
$wsALC.Range("I20").Value = 1294.4286
$wsALC.Range("K20").Value = 1294.4286
$wsALC.Range("M20").Value = -1064.4286

# ALC row 35 (G=1965)
$wsALC.Range("H35").Value = 1294.4286
$wsALC.Range("I35").Value = 1294.4286
$wsALC.Range("K35").Value = 1294.4286
$wsALC.Range("M35").Value = -915.4286

# ALC row 76 (G=12602)
$wsALC.Range("H76").Value = 5000
$wsALC.Range("J76").Value = 5000
$wsALC.Range("L76").Value = 5000
$wsALC.Range("N76").Value = -5630

# ALC row 79 (G=12602)
$wsALC.Range("H79").Value = 5000
$wsALC.Range("J79").Value = 5000
$wsALC.Range("L79").Value = 5000
$wsALC.Range("N79").Value = -7184

# ALC row 87 (G=10651)
$wsALC.Range("H87").Value = 31880.53
$wsALC.Range("J87").Value = 31880.53
$wsALC.Range("L87").Value = 31880.53
$wsALC.Range("N87").Value = -34376.53

# ALC row 90 (G=10651)
$wsALC.Range("H90").Value = 31880.53
$wsALC.Range("J90").Value = 31880.53
$wsALC.Range("L90").Value = 95641.59
$wsALC.Range("N90").Value = -108121.59

# ALC row 132 (G=44049)
$wsALC.Range("H132").Value = 3859.1333
$wsALC.Range("I132").Value = 2324.1667
$wsALC.Range("K132").Value = 6972.500100000001
$wsALC.Range("M132").Value = -4442.500100000001

# ALC row 135 (G=44047)
$wsALC.Range("H135").Value = 3692.5
$wsALC.Range("I135").Value = 3692.5
$wsALC.Range("K135").Value = 33232.5
$wsALC.Range("M135").Value = -30697.5

# ALC row 137 (G=44013)
$wsALC.Range("H137").Value = 2497.3333
$wsALC.Range("I137").Value = 2664
$wsALC.Range("J137").Value = 2330.6667
$wsALC.Range("K137").Value = 7992
$wsALC.Range("L137").Value = 6992.000100000001
$wsALC.Range("M137").Value = -5442
$wsALC.Range("N137").Value = -12092.0001

# ARM row 8 (G=3011)
$wsARM.Range("H8").Value = 25001
$wsARM.Range("I8").Value = 1005
$wsARM.Range("K8").Value = 1005
$wsARM.Range("M8").Value = -861

# ARM row 13 (G=2656)
$wsARM.Range("H13").Value = 9999
$wsARM.Range("I13").Value = 0
$wsARM.Range("J13").Value = 9999
$wsARM.Range("K13").Value = 0
$wsARM.Range("L13").Value = 9999
$wsARM.Range("M13").Value = ""
$wsARM.Range("N13").Value = -10287

# ARM row 45 (G=27714)
$wsARM.Range("H45").Value = 2599.8
$wsARM.Range("J45").Value = 2599.8
$wsARM.Range("L45").Value = 2599.8
$wsARM.Range("N45").Value = -3353.8

# ARM row 74 (G=44000)
$wsARM.Range("H74").Value = 1024.5
$wsARM.Range("I74").Value = 1024.5
$wsARM.Range("K74").Value = 1024.5
$wsARM.Range("M74").Value = -150.5

# ARM row 77 (G=44000)
$wsARM.Range("H77").Value = 1024.5
$wsARM.Range("I77").Value = 1024.5
$wsARM.Range("K77").Value = 5122.5
$wsARM.Range("M77").Value = -754.5

# ARM row 88 (G=12530)
$wsARM.Range("H88").Value = 1721.8462
$wsARM.Range("J88").Value = 1863.5714
$wsARM.Range("L88").Value = 1863.5714
$wsARM.Range("N88").Value = -2675.5714

# ARM row 91 (G=12530)
$wsARM.Range("H91").Value = 1721.8462
$wsARM.Range("J91").Value = 1863.5714
$wsARM.Range("L91").Value = 1863.5714
$wsARM.Range("N91").Value = -4671.5714

# ARM row 110 (G=27708)
$wsARM.Range("H110").Value = 164.25
$wsARM.Range("I110").Value = 163.57143
$wsARM.Range("K110").Value = 163.57143
$wsARM.Range("M110").Value = 1881.42857

# ARM row 122 (G=36168)
$wsARM.Range("H122").Value = 1687.4445
$wsARM.Range("I122").Value = 1648.375
$wsARM.Range("J122").Value = 2000
$wsARM.Range("K122").Value = 4945.125
$wsARM.Range("L122").Value = 6000
$wsARM.Range("M122").Value = -2495.125
$wsARM.Range("N122").Value = -10900

# BSM row 20 (G=14149)
$wsBSM.Range("H20").Value = 1242.4286
$wsBSM.Range("I20").Value = 949.5
$wsBSM.Range("K20").Value = 949.5
$wsBSM.Range("M20").Value = -702.5

# BSM row 99 (G=19943)
$wsBSM.Range("H99").Value = 1073.6471
$wsBSM.Range("I99").Value = 750.13336
$wsBSM.Range("K99").Value = 750.13336
$wsBSM.Range("M99").Value = 747.86664

# CUL row 4 (G=4650)
$wsCUL.Range("H4").Value = 55263304
$wsCUL.Range("I4").Value = 2777930.2
$wsCUL.Range("J4").Value = 1000000000
$wsCUL.Range("K4").Value = 8333790.600000001
$wsCUL.Range("L4").Value = 3000000000
$wsCUL.Range("M4").Value = -8333678.600000001
$wsCUL.Range("N4").Value = -3000000224

# CUL row 6 (G=4639)
$wsCUL.Range("H6").Value = 115.86667
$wsCUL.Range("I6").Value = 115.86667
$wsCUL.Range("K6").Value = 347.60001
$wsCUL.Range("M6").Value = -234.60001

# CUL row 45 (G=29501)
$wsCUL.Range("H45").Value = 5033
$wsCUL.Range("I45").Value = 0
$wsCUL.Range("J45").Value = 5033
$wsCUL.Range("K45").Value = 0
$wsCUL.Range("L45").Value = 15099
$wsCUL.Range("M45").Value = ""
$wsCUL.Range("N45").Value = -16163

# CUL row 50 (G=4725)
$wsCUL.Range("H50").Value = 1054.4546
$wsCUL.Range("I50").Value = 300
$wsCUL.Range("J50").Value = 1222.1111
$wsCUL.Range("K50").Value = 900
$wsCUL.Range("L50").Value = 3666.3333
$wsCUL.Range("M50").Value = -419
$wsCUL.Range("N50").Value = -4628.3333

# CUL row 53 (G=4725)
$wsCUL.Range("H53").Value = 1054.4546
$wsCUL.Range("I53").Value = 300
$wsCUL.Range("J53").Value = 1222.1111
$wsCUL.Range("K53").Value = 900
$wsCUL.Range("L53").Value = 3666.3333
$wsCUL.Range("M53").Value = -419
$wsCUL.Range("N53").Value = -4628.3333

# CUL row 80 (G=12890)
$wsCUL.Range("H80").Value = 3714
$wsCUL.Range("J80").Value = 3714
$wsCUL.Range("L80").Value = 11142
$wsCUL.Range("N80").Value = -13014

# CUL row 83 (G=12890)
$wsCUL.Range("H83").Value = 3714
$wsCUL.Range("J83").Value = 3714
$wsCUL.Range("L83").Value = 33426
$wsCUL.Range("N83").Value = -42786

# CUL row 133 (G=44073)
$wsCUL.Range("H133").Value = 0
$wsCUL.Range("I133").Value = 0
$wsCUL.Range("K133").Value = 0
$wsCUL.Range("M133").Value = ""

# LTW row 7 (G=36249)
$wsLTW.Range("H7").Value = 4301
$wsLTW.Range("I7").Value = 2666.6667
$wsLTW.Range("J7").Value = 6752.5
$wsLTW.Range("K7").Value = 2666.6667
$wsLTW.Range("L7").Value = 6752.5
$wsLTW.Range("M7").Value = -2554.6667
$wsLTW.Range("N7").Value = -6976.5

# LTW row 122 (G=36247)
$wsLTW.Range("H122").Value = 4656.4
$wsLTW.Range("I122").Value = 4424.6665
$wsLTW.Range("J122").Value = 5583.3335
$wsLTW.Range("K122").Value = 13273.9995
$wsLTW.Range("L122").Value = 16750.0005
$wsLTW.Range("M122").Value = -10823.9995
$wsLTW.Range("N122").Value = -21650.0005

# LTW row 126 (G=36249)
$wsLTW.Range("H126").Value = 4301
$wsLTW.Range("I126").Value = 2666.6667
$wsLTW.Range("J126").Value = 6752.5
$wsLTW.Range("K126").Value = 8000.000100000001
$wsLTW.Range("L126").Value = 20257.5
$wsLTW.Range("M126").Value = -5530.000100000001
$wsLTW.Range("N126").Value = -25197.5

# LTW row 136 (G=44060)
$wsLTW.Range("H136").Value = 2884.5
$wsLTW.Range("I136").Value = 2826.75
$wsLTW.Range("J136").Value = 3000
$wsLTW.Range("K136").Value = 8480.25
$wsLTW.Range("L136").Value = 9000
$wsLTW.Range("M136").Value = -5930.25
$wsLTW.Range("N136").Value = -14100

# WVR row 8 (G=2999)
$wsWVR.Range("H8").Value = 41666.668
$wsWVR.Range("I8").Value = 0
$wsWVR.Range("J8").Value = 41666.668
$wsWVR.Range("K8").Value = 0
$wsWVR.Range("L8").Value = 41666.668
$wsWVR.Range("M8").Value = ""
$wsWVR.Range("N8").Value = -41946.668

# WVR row 81 (G=12596)
$wsWVR.Range("H81").Value = 2000
$wsWVR.Range("I81").Value = 2000
$wsWVR.Range("J81").Value = 0
$wsWVR.Range("K81").Value = 4000
$wsWVR.Range("L81").Value = 0
$wsWVR.Range("M81").Value = -2939
$wsWVR.Range("N81").Value = ""

# WVR row 84 (G=12596)
$wsWVR.Range("H84").Value = 2000
$wsWVR.Range("I84").Value = 2000
$wsWVR.Range("J84").Value = 0
$wsWVR.Range("K84").Value = 20000
$wsWVR.Range("L84").Value = 0
$wsWVR.Range("M84").Value = -14696
$wsWVR.Range("N84").Value = ""
